# Append the latest EUR->ARS quote as a new row (row 8) at the bottom of the
# rate-history table, mirroring the existing rows (plain text cells for the
# date, time and quote columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date-looking / time-looking strings to be stored as literal text
# (matching the other rows) instead of letting Excel auto-convert them into
# date/time serial values. Briefly switching the number format to "Text"
# before the assignment achieves that; resetting the style back to "Normal"
# afterwards avoids leaving a stray number-format on the new cells.
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "2025-09-06"
$ws.Range("A8").Style = "Normal"

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "21:20:30"
$ws.Range("B8").Style = "Normal"

$ws.Range("C8").Value = "1.00 EUR = 1614.4992 ARS"
